# "Add files via upload" -- reset the Score sheet back to a blank template:
# drop the per-heat time/rep/score figures in columns C:G for every
# competitor row, leaving just the Team (A) and Name (B) labels, and
# leave the selection on M13 (matching the uploaded file's last cursor
# position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")
$ws.Activate()

# Clear out the Minute1/Second1/Rep1 (and the Rep1 sum formula) values
# for every data row (2-21), without touching A/B or disturbing the
# untouched H:O columns further along those same rows.
$ws.Range("C2:G21").ClearContents()

# Move the active selection to where the author last left it.
$ws.Range("M13").Select()
